# Apply cryptos list update (cell text values) to match target diff.
# Numeric-looking Price values need a leading apostrophe so Excel keeps
# them as text (matching the original inlineStr cell type) instead of
# auto-converting them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.161.54"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "3.464.20"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'585.45"
$ws.Range("E5").Value = "  +5.40%  "
$ws.Range("D6").Value = "'191.26"
$ws.Range("E6").Value = "  +9.47%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.460.69"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "'57.57"
$ws.Range("E12").Value = "  +7.39%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "'9.53"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").Value = "4.011.32"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "'19.01"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "3.462.32"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "67.147.42"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("E19").Value = "  +2.27%  "
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("E21").Value = "  +3.25%  "
$ws.Range("D22").Value = "'482.96"
$ws.Range("E22").Value = "  +6.26%  "
$ws.Range("D23").Value = "'5.34"
$ws.Range("E23").Value = "  +8.85%  "
$ws.Range("D24").Value = "'16.84"
$ws.Range("E24").Value = "  +17.24%  "
$ws.Range("E25").Value = "  +7.47%  "
$ws.Range("D26").Value = "'90.30"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").Value = "'3.00"
$ws.Range("E27").Value = "  +3.96%  "
$ws.Range("D28").Value = "'11.04"
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("D29").Value = "'9.09"
$ws.Range("E29").Value = "  +4.38%  "
$ws.Range("D30").Value = "'31.47"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +14.67%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'604.17"
$ws.Range("E32").Value = "  +4.61%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.87"
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").Value = "'64.47"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("E35").Value = "  +4.64%  "
$ws.Range("E36").Value = "  +4.89%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'37.60"
$ws.Range("E38").Value = "  +5.36%  "
$ws.Range("D39").Value = "'0.393"
$ws.Range("E39").Value = "  +5.71%  "
$ws.Range("D40").Value = "'3.51"
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").Value = "0.0₃0761"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("D42").Value = "3.225.62"
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("E43").Value = "  +7.11%  "
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("D45").Value = "'2.89"
$ws.Range("E45").Value = "  +28.14%  "
$ws.Range("D46").Value = "'2.60"
$ws.Range("E46").Value = "  +5.90%  "
$ws.Range("D47").Value = "'3.23"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "'0.136"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  +5.63%  "
$ws.Range("E50").Value = "  -0.03%  "
